$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)

$shape = $null
foreach ($sh in $s.Shapes) {
    if ($sh.Name -eq "TextBox 5") {
        $shape = $sh
    }
}

# Grow the textbox slightly (5796000 EMU -> 5832000 EMU, i.e. 456.37795... -> 459.21259... points)
$shape.Height = 459.21264

# Append a period to the final bullet's run text, keeping the rest of the
# paragraph/run structure (and the rest of the text body) untouched.
$tf = $shape.TextFrame
$tr = $tf.TextRange
$lastParagraph = $tr.Paragraphs(7, 1)
$lastRun = $lastParagraph.Runs(1, 1)
$lastRun.Text = "Passenger boarding in a queue increases the likelihood of choosing metro compared to crowded boarding."
